# "Roll Back Part -2"
# The underlying data edit captured by the diff is a rollback of the
# quantity/count value in Sheet1!A2 (150 -> 1000), together with moving the
# active selection to C5. (The remaining hunks in the source diff --
# fileVersion/rupBuild, xr:revisionPtr documentId, workbookView window
# geometry, and the x14ac:dyDescent / column-width jitter -- are artifacts
# of the authoring app's build/version re-saving the file and are not part
# of the addressable Excel object model, so they are not reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Roll back A2 from 150 to 1000.
$ws.Range("A2").Value = 1000

# Move/restore the active selection to C5, as shown in the diff.
$ws.Range("C5").Select()
